$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 432, pushing existing rows 432-447 down to 433-448.
$ws.Rows("432").Insert()

# Populate the newly inserted row 432 with a new price record (same station/category
# context as the surrounding rows, new cultivar/grade/price data).
$ws.Range("A432").Value = 7
$ws.Range("B432").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C432").Value = "Ñuble"
$ws.Range("D432").Value = 45267
$ws.Range("E432").Value = 16
$ws.Range("F432").Value = "Fruta"
$ws.Range("G432").Value = 100103
$ws.Range("H432").Value = "Frutos de hueso (carozo)"
$ws.Range("I432").Value = 100103004
$ws.Range("J432").Value = "Durazno"
$ws.Range("K432").Value = "Florida King"
$ws.Range("L432").Value = "Primera"
$ws.Range("M432").Value = 150
$ws.Range("N432").Value = 18000
$ws.Range("O432").Value = 18000
$ws.Range("P432").Value = 18000
$ws.Range("Q432").Value = "$/caja 15 kilos granel"
$ws.Range("R432").Value = "Región de O'Higgins"
$ws.Range("S432").Value = 1200
$ws.Range("T432").Value = 15
